# Bind the already-open workbook / active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header typo: "TRP_*" -> "TPR_*" for the 6 TPR columns (X1:AC1) ---
# Values keep the same column order, only the prefix spelling changes.
$ws.Range("X1").Value  = "TPR_6"
$ws.Range("Y1").Value  = "TPR_10"
$ws.Range("Z1").Value  = "TPR_15"
$ws.Range("AA1").Value = "TPR_18"
$ws.Range("AB1").Value = "TPR_6FFF"
$ws.Range("AC1").Value = "TPR_10FFF"

# --- Match the saved selection / view state from the edit: X1:AC1 selected, active cell X1 ---
[void]$ws.Range("X1:AC1").Select()
